$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains text formatting so numeric-looking strings
# (e.g. "1.00", "6.34") are not coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '63.679.48'
$ws.Range('D3').Value = '3.413.42'
$ws.Range('E3').Value = '  +2.41%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '569.57'
$ws.Range('E5').Value = '  +1.58%  '
$ws.Range('D6').Value = '157.05'
$ws.Range('E6').Value = '  +3.85%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '3.412.46'
$ws.Range('D9').Value = '0.544'
$ws.Range('E9').Value = '  +1.88%  '
$ws.Range('E10').Value = '  -0.15%  '
$ws.Range('E11').Value = '  +4.26%  '
$ws.Range('E12').Value = '  -0.24%  '
$ws.Range('D13').Value = '3.999.55'
$ws.Range('E13').Value = '  +2.23%  '
$ws.Range('E14').Value = '  -3.03%  '
$ws.Range('E15').Value = '  +8.74%  '
$ws.Range('D16').Value = '27.23'
$ws.Range('E16').Value = '  +2.04%  '
$ws.Range('D17').Value = '63.680.11'
$ws.Range('E17').Value = '  +2.05%  '
$ws.Range('D18').Value = '3.423.25'
$ws.Range('E18').Value = '  +3.80%  '
$ws.Range('E19').Value = '  -1.02%  '
$ws.Range('D20').Value = '14.06'
$ws.Range('E20').Value = '  +2.28%  '
$ws.Range('D21').Value = '378.54'
$ws.Range('E21').Value = '  -1.24%  '
$ws.Range('D22').Value = '8.06'
$ws.Range('E22').Value = '  -3.58%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('E24').Value = '  +2.98%  '
$ws.Range('E25').Value = '  -0.46%  '
$ws.Range('D26').Value = '0.0000120'
$ws.Range('E26').Value = '  +27.73%  '
$ws.Range('D27').Value = '9.37'
$ws.Range('E27').Value = '  +4.44%  '
$ws.Range('D28').Value = '0.179'
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('D30').Value = '6.03'
$ws.Range('E30').Value = '  +8.08%  '
$ws.Range('D31').Value = '1.36'
$ws.Range('E31').Value = '  +4.75%  '
$ws.Range('E32').Value = '  +1.74%  '
$ws.Range('D33').Value = '23.27'
$ws.Range('E33').Value = '  +2.00%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').Value = '6.33'
$ws.Range('E35').Value = '  -3.79%  '
$ws.Range('D36').Value = '6.77'
$ws.Range('E36').Value = '  +1.18%  '
$ws.Range('D37').Value = '159.14'
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('D38').Value = '1.44'
$ws.Range('E38').Value = '  -1.74%  '
$ws.Range('D39').Value = '2.984.48'
$ws.Range('E39').Value = '  +6.98%  '
$ws.Range('D40').Value = '0.0758'
$ws.Range('E40').Value = '  +3.07%  '
$ws.Range('D41').Value = '26.99'
$ws.Range('E41').Value = '  +0.63%  '
$ws.Range('D42').Value = '1.82'
$ws.Range('E42').Value = '  -3.33%  '
$ws.Range('D43').Value = '0.0316'
$ws.Range('E43').Value = '  +0.50%  '
$ws.Range('D44').Value = '41.97'
$ws.Range('E44').Value = '  +3.89%  '
$ws.Range('E45').Value = '  +2.87%  '
$ws.Range('D46').Value = '4.31'
$ws.Range('E46').Value = '  +1.50%  '
$ws.Range('E47').Value = '  +5.83%  '
$ws.Range('E48').Value = '  +3.81%  '
$ws.Range('D49').Value = '2.19'
$ws.Range('E49').Value = '  +23.04%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '6.34'
$ws.Range('E50').Value = '  +1.01%  '
$ws.Range('B51').Value = 'Bittensor'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D51').Value = '294.33'
$ws.Range('E51').Value = '  +2.59%  '
